$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'the hidden rainbow',
    'the hip',
    'the holy cow',
    'the iliac crest',
    'the key muscle of yoga',
    'the knee joint',
    'the limited black leggings',
    'the limited clothing women pants',
    'the marathon clothing',
    'the marathon stick',
    'the new ring of power',
    'the perfect sculpt leggings',
    'the record runner',
    'the scare squad',
    'the shins shirt',
    'the stick marathon',
    'the white woman on the green bicycle',
    'thermal compression',
    'thermal compression capris',
    'thermal compression women',
    'thermal for sports',
    'thermal leggings for women',
    'thermal running',
    'thermal running tights women',
    'thermal shorts',
    'thermal shorts for women',
    'thermal sports wear for men',
    'thick band leggings',
    'thick band leggings for women',
    'thick black leggings',
    'thick black womens tights',
    'thick leggings for women winter',
    'thick tights',
    'thick tights for women',
    'thick winter leggings for women',
    'thick winter tights',
    'thigh brace for hamstring',
    'thigh brace for women',
    'thigh compression brace',
    'thigh compression leggings',
    'thigh compression pants',
    'thigh compression shorts',
    'thigh compression support',
    'thigh compression tights',
    'thigh exercise',
    'thigh fitness',
    'thigh high',
    'thigh high 5x',
    'thigh high athletic',
    'thigh high black',
    'thigh high black tights',
    'thigh high glue',
    'thigh high leggings',
    'thigh high leggings winter',
    'thigh high rainbow',
    'thigh high tights',
    'thigh high tights for women',
    'thigh high workout leggings',
    'thigh highs',
    'thigh highs for women',
    'thigh highs with back seam',
    'thigh leggings for women',
    'thigh length compression',
    'thigh length leggings for women',
    'thigh pocket',
    'thigh ring',
    'thigh ring exercise',
    'thigh slimmer',
    'thigh slimmer men',
    'thigh slimmer shorts',
    'thigh slimmers for women',
    'thigh squeeze',
    'thigh support compression',
    'thigh tape',
    'thigh workout',
    'thigh workout equipment for women',
    'thighs exercise',
    'thighs for women',
    'thighs workout',
    'thights woman',
    'thin blue line wedding ring women',
    'thin leggings for summer',
    'thin tights',
    'thin tights for women',
    'think tank pro',
    'three quarter compression pants',
    'three quarter compression tights',
    'three stripe leggings',
    'thumb compression support',
    'thumb hole stock',
    'thumb joint',
    'thumb support compression',
    'tight and bright',
    'tight band',
    'tight basketball shorts',
    'tight black jacket women',
    'tight black pants',
    'tight black pants women',
    'tight black shorts women',
    'tight black tank'
)

for ($i = 1; $i -le $values.Count; $i++) {
    $ws.Cells.Item($i, 1).Value = $values[$i - 1]
}
